# refactoring + UML update !
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO_TEAM")
$ws.Activate()

# Insert a new row at row 5 (shifts existing rows 5..23 down to 6..24)
$ws.Rows.Item(5).Insert()

# Row 6 (was old row 5) - new description + changed priority
# (shared string for this text must be allocated before the row 5 text below)
$ws.Range("A6").Value = "IA ennemi?"
$ws.Range("C6").Value = "Bas"

# New task inserted at row 5
$ws.Range("A5").Value = "Initialiser la liste de Filtres (avec une Factory peut etre) (Gesture)"
$ws.Range("B5").Value = "ND"
$ws.Range("C5").Value = "Haut"
$ws.Range("D5").Value = "Ouvert"

# Row 7 (was old row 6) - new description + changed priority
$ws.Range("A7").Value = "Mettre a jour le diagramme de classe en cas de mofication du CODE !!!!!!"
$ws.Range("C7").Value = "Critique"

# All tasks are now marked as "Ouvert" instead of "Corrigé"
$ws.Range("D2:D23").Value = "Ouvert"

# Update selection to match authored state
$ws.Range("A6").Select()
